$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------------
# 1. Build the two new cell styles (cellXfs) using helper cells far outside
#    the used range, then propagate them with Copy/PasteSpecial(Formats) so
#    that no throw-away intermediate styles get recorded for the real cells.
# ---------------------------------------------------------------------------

# Style A: default number format, centered horizontally & vertically.
$helperA = $ws.Range("Z100")
$helperA.HorizontalAlignment = $xlCenter
$helperA.VerticalAlignment = $xlCenter
$helperA.Copy()
$ws.Range("A1:O6").PasteSpecial(-4122)

# Style B: custom date number format (yyyy-mm-dd;), centered horizontally &
# vertically - used only by column D.
$helperB = $ws.Range("Z101")
$helperB.NumberFormat = "yyyy-mm-dd;"
$helperB.HorizontalAlignment = $xlCenter
$helperB.VerticalAlignment = $xlCenter
$helperB.Copy()
$ws.Range("D1:D6").PasteSpecial(-4122)

$helperA.Clear()
$helperB.Clear()

# ---------------------------------------------------------------------------
# 2. Update header row text / add new header cells.
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

# ---------------------------------------------------------------------------
# 3. Fill in the previously-missing values (TotalConfirmedNewCases / G and
#    TotalNewDeaths / I) for the existing rows, and add the brand new row 6.
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("I2").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("I5").Value = 0

$ws.Range("A6").Value = 71
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = "LATIN AMER. & CARIB    "
$ws.Range("D6").Value = 43921
$ws.Range("E6").Value = "Anguilla"
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "Local transmission"
$ws.Range("K6").Value = 4
$ws.Range("L6").Value = 5384

# Columns M:O are intentionally left blank (no values) - they already exist
# in the sheet (and the dimension already reaches column O) purely because
# they were styled as part of the A1:O6 PasteSpecial(Formats) call above.

# ---------------------------------------------------------------------------
# 4. Column widths: 27 "characters" for columns A through O.
#    (ColumnWidth uses Excel's character-width units; the stored <col width>
#    value is ColumnWidth + ~0.8333, so we back that out to land on 27.)
# ---------------------------------------------------------------------------
$ws.Range("A1:O1").ColumnWidth = 26.1666666666667

